# Fruta / hortaliza, semanal
# A new weekly price-report entry was added as the new row 4 (all existing
# rows 4-96 shift down by one, to 5-97).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 4, pushing existing rows 4..96 down to 5..97.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new record.
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44643
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100109
$ws.Range("H4").Value = "Uva"
$ws.Range("I4").Value = 100109001
$ws.Range("J4").Value = "Uva"
$ws.Range("K4").Value = "Thompson seedless"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("Q4").Value = "`$/bandeja 18 kilos"
$ws.Range("R4").Value = "Región de Coquimbo"
$ws.Range("S4").Value = 972
$ws.Range("T4").Value = 18
